$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Grow the table from A1:G3 to A1:N3 by inserting new columns at the right
# spots. Inserting (rather than just overwriting trailing columns) makes the
# new header cells inherit the existing bold/bordered header style, and
# naturally shifts the existing "ACT"/"SAT"/"GPA" columns into their new
# homes (ACT Composite / SAT Combined / GPA).
# ---------------------------------------------------------------------------
$ws.Columns("D:E").Insert()   # room for RAI, Admit Score before Major
$ws.Columns("G:H").Insert()   # room for ACT Math, ACT English before (old) ACT
$ws.Columns("J:K").Insert()   # room for SAT Math, SAT Reading before (old) SAT
$ws.Columns("N:N").Insert()   # room for HS Percentile after GPA

# ---------------------------------------------------------------------------
# Header row (row 1):
#   A Name | B Total Amount | C Value | D RAI | E Admit Score | F Major |
#   G ACT Math | H ACT English | I ACT Composite |
#   J SAT Math | K SAT Reading | L SAT Combined | M GPA | N HS Percentile
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Total Amount"
$ws.Range("C1").Value = "Value"
$ws.Range("D1").Value = "RAI"
$ws.Range("E1").Value = "Admit Score"
$ws.Range("F1").Value = "Major"
$ws.Range("G1").Value = "ACT Math"
$ws.Range("H1").Value = "ACT English"
$ws.Range("I1").Value = "ACT Composite"
$ws.Range("J1").Value = "SAT Math"
$ws.Range("K1").Value = "SAT Reading"
$ws.Range("L1").Value = "SAT Combined"
$ws.Range("M1").Value = "GPA"
$ws.Range("N1").Value = "HS Percentile"

# ---------------------------------------------------------------------------
# Row 2 - "Test One"
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Test One"
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 8000
$ws.Range("D2").Value = 315
$ws.Range("E2").Value = 26
$ws.Range("F2").Value = "All"
$ws.Range("G2").Value = 25
$ws.Range("H2").Value = 27
$ws.Range("I2").Value = 26
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 400
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 96

# ---------------------------------------------------------------------------
# Row 3 - "Test Two" (was "Cool Kids Club")
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Test Two"
# Leading apostrophe forces these to stay text (matches source: "50"/"10000"
# stored as strings, not numbers).
$ws.Range("B3").Value = "'50"
$ws.Range("C3").Value = "'10000"
$ws.Range("D3").Value = 330
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = "All"
$ws.Range("G3").Value = 32
$ws.Range("H3").Value = 28
$ws.Range("I3").Value = 30
$ws.Range("J3").Value = 700
$ws.Range("K3").Value = 620
$ws.Range("L3").Value = 1320
$ws.Range("M3").Value = 3.8
$ws.Range("N3").Value = 95
